$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A24").Value = "JD_023"
$ws.Range("B24").Value = "Junior React Engineer"
$ws.Range("C24").Value = "We are seeking a Software Engineer to build and maintain high-quality software solutions.
Work with global teams to drive innovation and deliver scalable applications.
Join Akkodis and be part of a tech-driven, collaborative environment."
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 2

$ws.Rows(24).AutoFit()
